$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A59").Value = 'What''s the maximum number of lithology types in a log?'
$ws.Range("B59").Value = 'The maximum number of lithology types that can be represented in a log is 450.'
$ws.Range("A60").Value = 'What''s the maximum number of lithology types in a log?'
$ws.Range("B60").Value = 'The maximum number of lithology types that can be represented in a log is 450.'
$ws.Range("A61").Value = 'How many tracks can you define in one ODF?'
$ws.Range("B61").NumberFormat = "@"
$ws.Range("B61").Value = '20.'
$ws.Range("A62").Value = 'How many tracks can you define in one ODF?'
$ws.Range("B62").NumberFormat = "@"
$ws.Range("B62").Value = '20.'
$ws.Range("A63").Value = 'How many curve shades can I create?'
$ws.Range("B63").Value = 'According to the document, you can create 250 curve shades.'
$ws.Range("A64").Value = 'How many curve shades can I create?'
$ws.Range("B64").Value = 'According to the document, you can create 250 curve shades.'
$ws.Range("A65").Value = 'How many curves can I load in one go?'
$ws.Range("B65").Value = 'According to the GEO application documentation, you can load up to 5 data files to form one curve.'
$ws.Range("A66").Value = 'How many curves can I load in one go?'
$ws.Range("B66").Value = 'According to the GEO application documentation, you can load up to 5 data files to form one curve.'
$ws.Range("A67").Value = 'What is the maximum number of headers I can display in my log?'
$ws.Range("B67").Value = 'The maximum number of headers you can display in your log is not explicitly stated in the provided documentation. However, it does mention that on the Edit tab, selecting "Headers and Trailers" opens the Plot Header and Trailer Specifications dialog box, which lists various options for customization. It does not provide a specific limit for the number of headers.'
$ws.Range("A68").Value = 'What is the maximum number of headers I can display in my log?'
$ws.Range("B68").Value = 'The maximum number of headers you can display in your log is not explicitly stated in the provided documentation. However, it does mention that on the Edit tab, selecting "Headers and Trailers" opens the Plot Header and Trailer Specifications dialog box, which lists various options for customization. It does not provide a specific limit for the number of headers.'
$ws.Range("A69").Value = 'How many tables can I have in my log?'
$ws.Range("B69").Value = 'You can have up to 100 tables in a log.'
$ws.Range("A70").Value = 'How many tables can I have in my log?'
$ws.Range("B70").Value = 'You can have up to 100 tables in a log.'
$ws.Range("A71").Value = 'What''s the maximum number of characters in a single text entry?'
$ws.Range("B71").Value = 'The maximum number of characters in a single text entry is 32000.'
$ws.Range("A72").Value = 'What''s the maximum number of characters in a single text entry?'
$ws.Range("B72").Value = 'The maximum number of characters in a single text entry is 32000.'
$ws.Range("A73").Value = 'How many symbols can I have in the plot at any one time?'
$ws.Range("B73").Value = 'According to the document, you can have 10000 symbols per plot.'
$ws.Range("A74").Value = 'How many symbols can I have in the plot at any one time?'
$ws.Range("B74").Value = 'According to the document, you can have 10000 symbols per plot.'
$ws.Range("A75").Value = 'How many scales can I define?'
$ws.Range("B75").Value = 'According to the document, you can define up to 23 scales.'
$ws.Range("A76").Value = 'How many scales can I define?'
$ws.Range("B76").Value = 'According to the document, you can define up to 23 scales.'
$ws.Range("A77").Value = 'What is the maximum number of data files I can load?'
$ws.Range("B77").Value = 'The maximum number of data files you can load in one go is Unlimited.'
$ws.Range("A78").Value = 'What is the maximum number of data files I can load?'
$ws.Range("B78").Value = 'The maximum number of data files you can load in one go is Unlimited.'
